# Daily attendance processing - reorder "Recorded By" names in column G.
#
# For rows where "System" is listed first alongside a real user
# (dnasr281@gmail.com or backup@backdoor.com), move "System" so it comes
# after that user. A trailing lowercase "system" token (if present) stays
# last. Likewise, "admin@admin.com, dnasr281@gmail.com" becomes
# "dnasr281@gmail.com, admin@admin.com". Entries that don't contain one of
# these two priority users (e.g. "System, admin@admin.com" or a lone
# "System"/"dnasr281@gmail.com") are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

$oldValues = @(
    "System, backup@backdoor.com, system",
    "System, dnasr281@gmail.com",
    "System, backup@backdoor.com",
    "admin@admin.com, dnasr281@gmail.com"
)
$newValues = @(
    "backup@backdoor.com, System, system",
    "dnasr281@gmail.com, System",
    "backup@backdoor.com, System",
    "dnasr281@gmail.com, admin@admin.com"
)

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($val -ne $null) {
        for ($i = 0; $i -lt $oldValues.Length; $i++) {
            if ($val -eq $oldValues[$i]) {
                $cell.Value2 = $newValues[$i]
            }
        }
    }
}
